# TCR 2025 Results - update Record and Ranking sheets with Top-N group results
$wb = $excel.ActiveWorkbook

$wsRecord  = $wb.Worksheets.Item("Record")
$wsRanking = $wb.Worksheets.Item("Ranking")

# ---------------------------------------------------------------------------
# Record sheet: fill in the date + test-error values for the latest KNN k=3
# classification run (rows 51-58, one row per "Top N" group).
# ---------------------------------------------------------------------------
$wsRecord.Range("A51").Value = 45952

$wsRecord.Range("D51").Value = 0.4091
$wsRecord.Range("D52").Value = 0.3182
$wsRecord.Range("D53").Value = 0.3182
$wsRecord.Range("D54").Value = 0.2727
$wsRecord.Range("D55").Value = 0.2727
$wsRecord.Range("D56").Value = 0.3636
$wsRecord.Range("D57").Value = 0.3182
$wsRecord.Range("D58").Value = 0.3182

# Highlight the best (Top 5) result with vertical centering.
$wsRecord.Range("D51").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Ranking sheet: record which group each rank corresponds to, and its
# accuracy value.
# ---------------------------------------------------------------------------
$wsRanking.Range("C51").Value = "Top 20"
$wsRanking.Range("D51").Value = 0.2727
$wsRanking.Range("C51").HorizontalAlignment = -4108
$wsRanking.Range("D51").HorizontalAlignment = -4108
$wsRanking.Range("C51").Font.Size = 12
$wsRanking.Range("D51").Font.Size = 12

$wsRanking.Range("C52").Value = "Top 25"
$wsRanking.Range("D52").Value = 0.2727
$wsRanking.Range("C52").HorizontalAlignment = -4108
$wsRanking.Range("D52").HorizontalAlignment = -4108
$wsRanking.Range("C52").Font.Size = 12
$wsRanking.Range("D52").Font.Size = 12

$wsRanking.Range("C53").Value = "Top 30"
$wsRanking.Range("D53").Value = 0.3636
$wsRanking.Range("C53").HorizontalAlignment = -4108
$wsRanking.Range("D53").HorizontalAlignment = -4108
$wsRanking.Range("C53").Font.Size = 12
$wsRanking.Range("D53").Font.Size = 12

$wsRanking.Range("C54").Value = "Top 10"
$wsRanking.Range("D54").Value = 0.3182
$wsRanking.Range("C54").HorizontalAlignment = -4108
$wsRanking.Range("D54").HorizontalAlignment = -4108
$wsRanking.Range("C54").Font.Size = 12
$wsRanking.Range("D54").Font.Size = 12

$wsRanking.Range("C55").Value = "Top 15"
$wsRanking.Range("D55").Value = 0.3182
$wsRanking.Range("C55").HorizontalAlignment = -4108
$wsRanking.Range("D55").HorizontalAlignment = -4108
$wsRanking.Range("C55").Font.Size = 12
$wsRanking.Range("D55").Font.Size = 12

$wsRanking.Range("C56").Value = "Top 35"
$wsRanking.Range("D56").Value = 0.3182
$wsRanking.Range("C56").HorizontalAlignment = -4108
$wsRanking.Range("D56").HorizontalAlignment = -4108
$wsRanking.Range("C56").Font.Size = 12
$wsRanking.Range("D56").Font.Size = 12

$wsRanking.Range("C57").Value = "Top 37"
$wsRanking.Range("D57").Value = 0.3182
$wsRanking.Range("C57").HorizontalAlignment = -4108
$wsRanking.Range("D57").HorizontalAlignment = -4108
$wsRanking.Range("C57").Font.Size = 12
$wsRanking.Range("D57").Font.Size = 12

$wsRanking.Range("C58").Value = "Top 5"
$wsRanking.Range("D58").Value = 0.4091
$wsRanking.Range("C58").HorizontalAlignment = -4108
$wsRanking.Range("C58").Font.Size = 12

# Highlight the best (Top 5) result with vertical centering.
$wsRanking.Range("D58").HorizontalAlignment = -4108
$wsRanking.Range("D58").Font.Size = 12
$wsRanking.Range("D58").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Selections: leave Record's selection over the block that was edited, and
# finish with Ranking as the active sheet/selection (matches the saved file).
# ---------------------------------------------------------------------------
$wsRecord.Range("C51:D58").Select()
$wsRanking.Range("C56:D57").Select()
